$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are unambiguous text (contain non-numeric characters)
# can be assigned directly; Excel will keep them as text automatically.
$ws.Range("D2").Value = "58.708.48"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "2.302.83"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("E6").Value = "  +1.15%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +2.79%  "
$ws.Range("D9").Value = "2.297.80"
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("E10").Value = "  -0.67%  "
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("E13").Value = "  -0.07%  "
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("D15").Value = "2.711.56"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("D16").Value = "58.537.55"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").Value = "2.282.61"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("E20").Value = "  -2.08%  "
$ws.Range("E21").Value = "  +0.58%  "
$ws.Range("E22").Value = "  +2.04%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("E27").Value = "  -1.22%  "
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("E30").Value = "  -1.85%  "
$ws.Range("D31").Value = "0.0₃0732"
$ws.Range("E31").Value = "  +0.42%  "
$ws.Range("E32").Value = "  +2.70%  "
$ws.Range("E33").Value = "  +1.02%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("E35").Value = "  +0.40%  "
$ws.Range("B36").Value = "USDe"
$ws.Range("C36").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  +0.18%  "
$ws.Range("E39").Value = "  +0.35%  "
$ws.Range("E40").Value = "  -0.60%  "
$ws.Range("E41").Value = "  -1.70%  "
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("E46").Value = "  +0.43%  "
$ws.Range("E47").Value = "  -1.07%  "
$ws.Range("E48").Value = "  -1.63%  "
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("B51").Value = "ZEEBU"
$ws.Range("C51").Value = "https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu"
$ws.Range("E51").Value = "  +0.05%  "

# Cells whose new values look like plain numbers must be forced to remain
# text (matching the original inlineStr cell type) by temporarily applying
# a text number format, then restoring the default style afterwards.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.335"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000134"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "317.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.384"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.94"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "294.13"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "141.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0955"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0499"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.557"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.95"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.62"
$ws.Range("D51").Style = "Normal"
